$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows (863-867) for a new weekly price report (date 2022-02-03,
# serial 44595), pushing the existing rows 863-877 down to 868-882.
$ws.Range("A863:A867").EntireRow.Insert()

# Common (constant) column values shared by every record in this block.
$colA = 6
$colB = "Mercado Mayorista Lo Valledor de Santiago"
$colC = "Metropolitana"
$colE = 13
$colF = 100112028
$colG = "Sandia"
$colH = "Sin especificar"
$colN = "`$/unidad"
$colO = "Región de O'Higgins"
$colQ = 1
$colR = "Hortaliza"

# New weekly data rows, in order, with the variety-specific values.
$newRows = @(
    @{ Row = 863; I = "Extra";    J = 3500; K = 2500; L = 2900; M = 2671 },
    @{ Row = 864; I = "Primera";  J = 6500; K = 2000; L = 2300; M = 2134 },
    @{ Row = 865; I = "Segunda";  J = 5300; K = 1500; L = 1800; M = 1636 },
    @{ Row = 866; I = "Super";    J = 3900; K = 3000; L = 3300; M = 3131 },
    @{ Row = 867; I = "Tercera";  J = 3100; K = 900;  L = 1300; M = 1068 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $colA
    $ws.Cells.Item($row, 2).Value = $colB
    $ws.Cells.Item($row, 3).Value = $colC
    $ws.Cells.Item($row, 4).Value = 44595
    $ws.Cells.Item($row, 5).Value = $colE
    $ws.Cells.Item($row, 6).Value = $colF
    $ws.Cells.Item($row, 7).Value = $colG
    $ws.Cells.Item($row, 8).Value = $colH
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $colN
    $ws.Cells.Item($row, 15).Value = $colO
    $ws.Cells.Item($row, 16).Value = $r.M
    $ws.Cells.Item($row, 17).Value = $colQ
    $ws.Cells.Item($row, 18).Value = $colR
}
